# Begin of big refont - rework the "pyinstaller" command paragraph:
#  1. Duplicate "Plannificateur>" just before "pyinstaller" (typo-ish repeat, as in the
#     authored edit) inside the command-line paragraph.
#  2. Merge the trailing `-import="h5py._proxy"` / ` ` / `Plannificateur.py` runs back
#     together cleanly (no stray autocorrect / smart-quote substitution).
#  3. Fold that whole paragraph into the next one (the "Use spyder working_exe as env"
#     paragraph), i.e. remove the paragraph break that used to separate them, while
#     keeping the following paragraph's formatting (language = en-US).

$d = $word.ActiveDocument

# --- Step 1: duplicate "Plannificateur>" right before "pyinstaller" ------------------
$find = $d.Content
$found = $find.Find.Execute("Plannificateur>pyinstaller", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $matchStart = $find.Start
    $insPos = $matchStart + 15   # length of "Plannificateur>" -> right before "pyinstaller"
    $delRange = $d.Range($insPos - 1, $insPos)
    $delRange.Delete()
    $newPos = $insPos - 1
    $insRange = $d.Range($newPos, $newPos)
    $insRange.InsertBefore(">Plannificateur>")
}

# --- Step 2: merge the closing quote / space / filename runs -------------------------
$find2 = $d.Content
$found2 = $find2.Find.Execute("Plannificateur.py", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $pyStart = $find2.Start
    $delRange2 = $d.Range($pyStart - 2, $pyStart)
    $delRange2.Delete()
    $newPos2 = $pyStart - 2
    $insRange2 = $d.Range($newPos2, $newPos2)
    $insRange2.InsertBefore([char]34 + " ")
}

# --- Step 3: merge the command-line paragraph into the following paragraph -----------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Plannificateur.py*") {
        $r = $p.Range
        $markRange = $d.Range($r.End - 1, $r.End)
        $markRange.Delete()
        break
    }
}
